# Code Plan.xlsx update -- "updates through 1/2/24 17:15 - all RecordMsg methods validated"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$nl = [char]10

# ---------------------------------------------------------------------------
# 1. Update the docstring of __init__ (row 2, column F) -- text grows with
#    ".df_errs from Excel file" appended to the second line.
# ---------------------------------------------------------------------------
$ws.Range("F2").Value = "Initialize all attributes" + $nl + "Import Error Codes to .df_errs from Excel file"

# ---------------------------------------------------------------------------
# 2. Insert three new method rows for the ErrorHandle class (RecordMsg,
#    GetBaseErrCode, SetReportErrCode) above the existing is_fail row, then
#    trim four trailing blank rows so the used range stays at H16 (matches
#    the target dimension A1:H16).
# ---------------------------------------------------------------------------
$ws.Rows("3:5").Insert()
$ws.Rows("17:20").Delete()

# --- Row 3: RecordMsg (procedure) -----------------------------------------
$ws.Range("A3").Value = "ErrorHandle"
$ws.Range("B3").Value = "RecordMsg (procedure)"
$ws.Range("G3").Value = "Process an error or warning into a msg string"

# --- Row 4: GetBaseErrCode (tall row because of the long docstring) -------
$ws.Range("A4").Value = "ErrorHandle"
$ws.Range("B4").Value = "GetBaseErrCode"
$ws.Range("F4").Value = 'Set iCodeBase to global variable iErrNotFound if no rows match .Locn in the .df_errs "Locn" column' + $nl + 'Set base_row to the .df_errs rows matching .Locn in the "Locn" column and matching "Base" in the Msg_String column' + $nl + 'if base_row is not empty, set .iCodeBase to the .df_errs "iCode" column value for base_row[0]'
$ws.Range("G4").Value = "Look up Base .df_errs code for .Locn"
$ws.Rows("4:4").RowHeight = 70.5

# --- Row 5: SetReportErrCode ------------------------------------------------
$ws.Range("A5").Value = "ErrorHandle"
$ws.Range("B5").Value = "SetReportErrCode"
$ws.Range("G5").Value = "Process an error or warning into a msg string"

# ---------------------------------------------------------------------------
# 3. Formatting: column A on the three new rows uses a lighter "continuation"
#    fill (Accent1 font/border, but a 40%-tinted Accent1 fill) while columns
#    B, C, D-G, H reuse the same look as the surrounding method rows.
#    NOTE: PasteSpecial only honours the first area of a multi-area
#    (comma-separated) destination range, so every destination is pasted
#    individually.
# ---------------------------------------------------------------------------
foreach ($r in 3, 4, 5) {
    $ws.Range("B2").Copy()
    $ws.Range("B$r").PasteSpecial(-4122)

    $ws.Range("C2").Copy()
    $ws.Range("C$r").PasteSpecial(-4122)

    $ws.Range("D2").Copy()
    $ws.Range("D$r").PasteSpecial(-4122)
    $ws.Range("E2").Copy()
    $ws.Range("E$r").PasteSpecial(-4122)
    $ws.Range("F2").Copy()
    $ws.Range("F$r").PasteSpecial(-4122)
    $ws.Range("G2").Copy()
    $ws.Range("G$r").PasteSpecial(-4122)

    $ws.Range("H2").Copy()
    $ws.Range("H$r").PasteSpecial(-4122)

    # Column A: start from the existing "ErrorHandle" accent style (keeps
    # the thin border + Accent1 cell style), then lighten the fill to the
    # new tint.
    $ws.Range("A2").Copy()
    $ws.Range("A$r").PasteSpecial(-4122)
    $ws.Range("A$r").Interior.ThemeColor = 5
}

$excel.CutCopyMode = 0
